$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'244.31"
$ws.Range("D3").Formula = "'23.78"
$ws.Range("D4").Formula = "'5.195"
$ws.Range("D5").Formula = "'0.05750"
$ws.Range("D6").Formula = "'6.453"
$ws.Range("D7").Formula = "'3.253"
$ws.Range("D8").Formula = "'0.8095"
$ws.Range("D9").Formula = "'0.8688"
$ws.Range("D10").Formula = "'0.1374"
$ws.Range("D11").Formula = "'0.06947"
$ws.Range("D12").Formula = "'0.03187"
$ws.Range("D13").Formula = "'0.03029"
$ws.Range("D14").Formula = "'0.09323"
$ws.Range("D15").Formula = "'3.815"
$ws.Range("D16").Formula = "'0.001525"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Formula = "'0.006182"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").Formula = "'0.001234"
$ws.Range("E19").Value = "18BitKanKAN"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Formula = "'0.004084"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Formula = "'0.00008700"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Formula = "'3.553"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Formula = "'2.146"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Formula = "'0.01013"
$ws.Range("E24").Value = "23OneONEBestin24h"
$ws.Range("D26").Formula = "'0.1296"
$ws.Range("D27").Formula = "'0.0002328"
$ws.Range("D40").Formula = "'0.03706"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Formula = "'0.006260"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Formula = "'0.1046"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Formula = "'0.002449"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("D44").Formula = "'0.007067"
$ws.Range("D45").Formula = "'0.00005294"
$ws.Range("D46").Formula = "'0.00000000750"
$ws.Range("D47").Formula = "'0.5299"
$ws.Range("D48").Formula = "'0.002064"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Formula = "'0.00002100"
$ws.Range("D50").Formula = "'0.0002000"
